$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on Price/Volume columns so numeric-looking strings stay text
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '67.075.88'
$ws.Range("E2").Value = '  -0.50%  '
$ws.Range("D3").Value = '2.468.83'
$ws.Range("E3").Value = '  -0.56%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '581.58'
$ws.Range("E5").Value = '  -1.04%  '
$ws.Range("D6").Value = '167.91'
$ws.Range("E6").Value = '  -2.03%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '0.515'
$ws.Range("E8").Value = '  -1.52%  '
$ws.Range("D9").Value = '2.467.11'
$ws.Range("E9").Value = '  -0.54%  '
$ws.Range("D10").Value = '0.135'
$ws.Range("E10").Value = '  -1.84%  '
$ws.Range("E11").Value = '  -0.38%  '
$ws.Range("D12").Value = '4.99'
$ws.Range("E12").Value = '  -1.98%  '
$ws.Range("E13").Value = '  -1.79%  '
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = '25.49'
$ws.Range("E14").Value = '  -2.59%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '2.917.62'
$ws.Range("E15").Value = '  -0.92%  '
$ws.Range("D16").Value = '66.918.53'
$ws.Range("E16").Value = '  -0.50%  '
$ws.Range("D17").Value = '0.0000170'
$ws.Range("E17").Value = '  -3.21%  '
$ws.Range("D18").Value = '2.479.58'
$ws.Range("E18").Value = '  -0.44%  '
$ws.Range("D19").Value = '11.13'
$ws.Range("E19").Value = '  -4.23%  '
$ws.Range("D20").Value = '7.56'
$ws.Range("E20").Value = '  -3.13%  '
$ws.Range("D21").Value = '353.63'
$ws.Range("E21").Value = '  -3.11%  '
$ws.Range("D22").Value = '4.02'
$ws.Range("E22").Value = '  -2.44%  '
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("D24").Value = '69.03'
$ws.Range("E24").Value = '  -2.88%  '
$ws.Range("D25").Value = '4.24'
$ws.Range("E25").Value = '  -6.70%  '
$ws.Range("E26").Value = '  -6.54%  '
$ws.Range("D27").Value = '9.16'
$ws.Range("E27").Value = '  -7.51%  '
$ws.Range("E28").Value = '  -0.25%  '
$ws.Range("D29").Value = '2.591.82'
$ws.Range("E29").Value = '  -0.39%  '
$ws.Range("D30").Value = '0.0₃0907'
$ws.Range("E30").Value = '  -4.85%  '
$ws.Range("D31").Value = '514.38'
$ws.Range("E31").Value = '  -2.56%  '
$ws.Range("E32").Value = '  -6.25%  '
$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").Value = '1.24'
$ws.Range("E33").Value = '  -4.63%  '
$ws.Range("B34").Value = 'PancakeSwap'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D34").Value = '1.77'
$ws.Range("E34").Value = '  -4.14%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("D36").Value = '0.119'
$ws.Range("E36").Value = '  -6.91%  '
$ws.Range("D37").Value = '158.02'
$ws.Range("E37").Value = '  -0.87%  '
$ws.Range("D38").Value = '18.63'
$ws.Range("E38").Value = '  +0.08%  '
$ws.Range("D39").Value = '18.40'
$ws.Range("E39").Value = '  -1.47%  '
$ws.Range("E40").Value = '  -4.34%  '
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("D42").Value = '1.67'
$ws.Range("E42").Value = '  -5.35%  '
$ws.Range("D43").Value = '0.327'
$ws.Range("E43").Value = '  -6.01%  '
$ws.Range("D44").Value = '4.79'
$ws.Range("E44").Value = '  -5.32%  '
$ws.Range("D45").Value = '2.34'
$ws.Range("E45").Value = '  -4.15%  '
$ws.Range("D46").Value = '38.65'
$ws.Range("E46").Value = '  -2.45%  '
$ws.Range("D47").Value = '140.80'
$ws.Range("E47").Value = '  -2.62%  '
$ws.Range("E48").Value = '  -5.40%  '
$ws.Range("D49").Value = '0.514'
$ws.Range("E49").Value = '  -5.46%  '
$ws.Range("B50").Value = 'Optimism'
$ws.Range("C50").Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range("D50").Value = '1.59'
$ws.Range("E50").Value = '  -4.84%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₆0252'
$ws.Range("E51").Value = '  -9.07%  '
